$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 62. All existing rows from 62 downward
# (62..86) shift down by one (to 63..87), which matches the diff: every
# row's data now equals the row above it from the previous layout, and a
# brand-new weekly observation is inserted at row 62.
$ws.Rows.Item(62).Insert()

# Populate the newly inserted row 62 with the new weekly record.
$ws.Range("A62").Value2 = 6
$ws.Range("B62").Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C62").Value2 = "Metropolitana"
$ws.Range("D62").Value2 = 44839
$ws.Range("E62").Value2 = 13
$ws.Range("F62").Value2 = "Fruta"
$ws.Range("G62").Value2 = 100108
$ws.Range("H62").Value2 = "Tropicales y subtropicales"
$ws.Range("I62").Value2 = 100108007
$ws.Range("J62").Value2 = "Coco"
$ws.Range("K62").Value2 = "Sin especificar"
$ws.Range("L62").Value2 = "Primera"
$ws.Range("M62").Value2 = 150
$ws.Range("N62").Value2 = 29000
$ws.Range("O62").Value2 = 30000
$ws.Range("P62").Value2 = 29500
$ws.Range("Q62").Value2 = "$/malla 20 unidades"
$ws.Range("R62").Value2 = "Perú"
$ws.Range("S62").Value2 = 1475
$ws.Range("T62").Value2 = 20
